$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 3 (section "2. Editar Tarea" - "Para que" run): split the single run
# " pueda actualizar detalles como el título, la descripción, las subtareas,
#  la fecha de finalización o la categoría, siempre que la tarea no esté
#  completada." into three runs, dropping the "subtareas" / "categoría"
# mentions. Do this first while absolute character offsets are still
# pristine (nothing earlier in the doc has shifted yet).
# ---------------------------------------------------------------------------

# Delete ' o la categoría,' (rightmost edit first so earlier offsets stay valid)
$rA = $d.Range(372, 388)
if ($rA.Text -ne " o la categoría,") {
    throw "Unexpected text at rA: [$($rA.Text)]"
}
$rA.Delete()

# Delete 's subtareas, la' (leaves 'la' from the preceding word + ' fecha' after)
$rB = $d.Range(335, 350)
if ($rB.Text -ne "s subtareas, la") {
    throw "Unexpected text at rB: [$($rB.Text)]"
}
$rB.Delete()

# Force the middle segment "descripción, la" to serialize as its own run by
# genuinely toggling a character property (set-then-restore), mirroring how
# Word splits runs when a user re-touches a sub-selection.
$mid = $d.Range(320, 335)
if ($mid.Text -ne "descripción, la") {
    throw "Unexpected text at mid: [$($mid.Text)]"
}
$mid.Bold = 1
$mid.Bold = 0

# ---------------------------------------------------------------------------
# Change 2 (section "1. Crear Tarea" - "Para que" run): drop ", subtareas"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("con título, descripción, subtareas, fecha de finalización y categoría.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "con título, descripción, fecha de finalización y categoría.", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 4: remove user stories 6-9 ("6. Ver Estadísticas de Productividad"
# through "9. Manejar Tareas Recurrentes") plus the trailing empty paragraph.
# ---------------------------------------------------------------------------
$startPara = $d.Paragraphs(21)
if ($startPara.Range.Text.TrimEnd([char]13) -ne "6. Ver Estadísticas de Productividad") {
    throw "Unexpected paragraph 21: [$($startPara.Range.Text)]"
}
$endPara = $d.Paragraphs(37)
if ($endPara.Range.Text.TrimEnd([char]13) -ne "") {
    throw "Unexpected paragraph 37: [$($endPara.Range.Text)]"
}
$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------------
# Change 1: insert a new bold "HISTORIAS DE USUARIO" paragraph before the
# first existing paragraph ("1. Crear Tarea"). Do this last so it doesn't
# shift the paragraph indices / character offsets used above.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
if ($firstPara.Range.Text.TrimEnd([char]13) -ne "1. Crear Tarea") {
    throw "Unexpected paragraph 1: [$($firstPara.Range.Text)]"
}
$firstPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs(1)
$newPara.Range.InsertBefore("HISTORIAS DE USUARIO")
